$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1279
$ws.Range("F7").Value = 837
$ws.Range("F8").Value = 34
$ws.Range("F9").Value = 6722
$ws.Range("F11").Value = 97
$ws.Range("F12").Value = 135
$ws.Range("F13").Value = 6384
$ws.Range("F14").Value = 121
$ws.Range("F15").Value = 266
$ws.Range("F19").Value = 4260
$ws.Range("F20").Value = 217
$ws.Range("F21").Value = 226
$ws.Range("F23").Value = 306
$ws.Range("F32").Value = 7761
$ws.Range("F34").Value = 1306
$ws.Range("F35").Value = 640
$ws.Range("F36").Value = 14
$ws.Range("F37").Value = 115
$ws.Range("F40").Value = 1544
$ws.Range("F41").Value = 202
$ws.Range("F42").Value = 885
$ws.Range("F43").Value = 38
$ws.Range("F44").Value = 3852
$ws.Range("F46").Value = 21
$ws.Range("F47").Value = 105
$ws.Range("F48").Value = 824

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 71
$ws.Range("F8").Value = 11
$ws.Range("F11").Value = 154
$ws.Range("F16").Value = 75

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 71
$ws.Range("F7").Value = 1279
$ws.Range("F11").Value = 837
$ws.Range("F12").Value = 34
$ws.Range("F13").Value = 6722
$ws.Range("F15").Value = 97
$ws.Range("F16").Value = 135
$ws.Range("F17").Value = 6384
$ws.Range("F18").Value = 121
$ws.Range("F19").Value = 266
$ws.Range("F21").Value = 4260
$ws.Range("F22").Value = 217
$ws.Range("F23").Value = 226
$ws.Range("F25").Value = 306
$ws.Range("F30").Value = 154
$ws.Range("F31").Value = 7761
$ws.Range("F33").Value = 1306
$ws.Range("F34").Value = 640
$ws.Range("F35").Value = 14
$ws.Range("F36").Value = 115
$ws.Range("F39").Value = 1544
$ws.Range("F40").Value = 202
$ws.Range("F41").Value = 885
$ws.Range("F42").Value = 38
$ws.Range("F43").Value = 3852
$ws.Range("F45").Value = 21
$ws.Range("F47").Value = 824
